$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- sheet2 ("description"): add the two new helper columns E:F ---

# Column widths for the new columns (engine quantizes to nearest 1/6 char;
# 24.8 lands closest to the target 25.7265625)
$ws2.Range("E1").ColumnWidth = 24.8
$ws2.Range("F1").ColumnWidth = 24.8

# Row 1: merged header cell E1:F1, centered horizontally, no content
$ws2.Range("E1").HorizontalAlignment = -4108
$ws2.Range("F1").HorizontalAlignment = -4108

# Row 2: E2 = "Neutral" style note (can be blank), centered
$e2 = $ws2.Range("E2")
$e2.Value = "เป็นค่าว่างได้"
$e2.Style = "Neutral"
$e2.HorizontalAlignment = -4108
$e2.VerticalAlignment = -4108

# Row 2: F2 = reuse the existing "Bad" centered style (already used by A2:C2)
$f2 = $ws2.Range("F2")
$f2.Value = "เพิ่มข้อมูล ต้องไม่เป็นค่าว่าง"
$ws2.Range("C2").Copy()
$f2.PasteSpecial(-4122)

# Merge the header cells after both retain their own formatting
$ws2.Range("E1:F1").Merge()

# Rows 3-5: touch E/F so empty cells materialize in the row, matching the
# surrounding row's default style, then drop any value back out
$ws2.Range("E3").Value = "x"
$ws2.Range("E3").ClearContents()
$ws2.Range("F3").Value = "x"
$ws2.Range("F3").ClearContents()

$ws2.Range("E4").Value = "x"
$ws2.Range("E4").ClearContents()
$ws2.Range("F4").Value = "x"
$ws2.Range("F4").ClearContents()

$ws2.Range("E5").Value = "x"
$ws2.Range("E5").ClearContents()
$ws2.Range("F5").Value = "x"
$ws2.Range("F5").ClearContents()

# --- tab/selection state: "description" sheet becomes the active tab ---
$ws1.Range("A6").Select()
$ws2.Activate()
$ws2.Range("C15").Select()

Write-Host "done"
